$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0.3299202878440668
$ws.Cells.Item(2, 2).Value = 0.3829666873843514
$ws.Cells.Item(2, 3).Value = 0.2822597518622293
$ws.Cells.Item(2, 4).Value = 0.3548941068322122
$ws.Cells.Item(2, 5).Value = 0.3038515979286224
$ws.Cells.Item(2, 11).Value = 2.309442014908468
$ws.Cells.Item(2, 12).Value = 2.68076681169046
$ws.Cells.Item(2, 13).Value = 1.975818263035606
$ws.Cells.Item(2, 14).Value = 2.484258747825486
$ws.Cells.Item(2, 15).Value = 2.126961185500357
$ws.Cells.Item(2, 16).Value = 63.77354
$ws.Cells.Item(2, 17).Value = 81.75593124946359
$ws.Cells.Item(2, 18).Value = 50.26230442364822
$ws.Cells.Item(2, 19).Value = 71.46314943887475
$ws.Cells.Item(2, 20).Value = 55.56684899504945
$ws.Cells.Item(2, 21).Value = 0.2036050555118729
$ws.Cells.Item(2, 22).Value = 0.2591105867630624
$ws.Cells.Item(2, 23).Value = 0.1490971507121034
$ws.Cells.Item(2, 24).Value = 0.232234898062347
$ws.Cells.Item(2, 25).Value = 0.1737161120343146
$ws.Cells.Item(2, 26).Value = 0.8600071441520277
$ws.Cells.Item(2, 27).Value = 0.915263696753041
$ws.Cells.Item(2, 28).Value = 0.7913138442740167
$ws.Cells.Item(2, 29).Value = 0.8902884558500727
$ws.Cells.Item(2, 30).Value = 0.8249862327812542
$ws.Cells.Item(3, 6).Value = 6.997067711031733
$ws.Cells.Item(3, 7).Value = 8.828010335236128
$ws.Cells.Item(3, 8).Value = 5.451119114364335
$ws.Cells.Item(3, 9).Value = 7.833502622234264
$ws.Cells.Item(3, 10).Value = 6.123672032234734
$ws.Cells.Item(3, 11).Value = 2.309032344640472
$ws.Cells.Item(3, 12).Value = 2.913243410627922
$ws.Cells.Item(3, 13).Value = 1.79886930774023
$ws.Cells.Item(3, 14).Value = 2.585055865337307
$ws.Cells.Item(3, 15).Value = 2.020811770637462
$ws.Cells.Item(3, 16).Value = 63.76815
$ws.Cells.Item(3, 17).Value = 75.13540296383746
$ws.Cells.Item(3, 18).Value = 55.97465215521801
$ws.Cells.Item(3, 19).Value = 68.77352711372066
$ws.Cells.Item(3, 20).Value = 59.30487705283106
$ws.Cells.Item(3, 21).Value = 0.2022668229566483
$ws.Cells.Item(3, 22).Value = 0.289500480703474
$ws.Cells.Item(3, 23).Value = 0.1176270291088261
$ws.Cells.Item(3, 24).Value = 0.2473911905271854
$ws.Cells.Item(3, 25).Value = 0.1554198925436388
$ws.Cells.Item(3, 26).Value = 0.8537973356140834
$ws.Cells.Item(3, 27).Value = 0.9381132992451295
$ws.Cells.Item(3, 28).Value = 0.7353486350461407
$ws.Cells.Item(3, 29).Value = 0.9012536280553676
$ws.Cells.Item(3, 30).Value = 0.7958025402878478
$ws.Cells.Item(4, 1).Value = 0.3299709879228907
$ws.Cells.Item(4, 2).Value = 0.3835433800069381
$ws.Cells.Item(4, 3).Value = 0.282490288356253
$ws.Cells.Item(4, 4).Value = 0.3549708526407032
$ws.Cells.Item(4, 5).Value = 0.3039131997519875
$ws.Cells.Item(4, 6).Value = 7.002209084319542
$ws.Cells.Item(4, 7).Value = 8.842692131684961
$ws.Cells.Item(4, 8).Value = 5.451334083599908
$ws.Cells.Item(4, 9).Value = 7.833899822736148
$ws.Cells.Item(4, 10).Value = 6.132871319663241
$ws.Cells.Item(4, 11).Value = 2.310539526478212
$ws.Cells.Item(4, 12).Value = 3.041834367661179
$ws.Cells.Item(4, 13).Value = 1.719625473759511
$ws.Cells.Item(4, 14).Value = 2.635556163608213
$ws.Cells.Item(4, 15).Value = 1.971112809645559
$ws.Cells.Item(4, 16).Value = 64.82935000000001
$ws.Cells.Item(4, 17).Value = 89.97420461149105
$ws.Cells.Item(4, 18).Value = 49.15547311527516
$ws.Cells.Item(4, 19).Value = 74.05936458671786
$ws.Cells.Item(4, 20).Value = 55.43313974299544
$ws.Cells.Item(4, 21).Value = 0.2016719693897677
$ws.Cells.Item(4, 22).Value = 0.30442753138228
$ws.Cells.Item(4, 23).Value = 0.1029117569653542
$ws.Cells.Item(4, 24).Value = 0.255087482693251
$ws.Cells.Item(4, 25).Value = 0.1464663332311515
$ws.Cells.Item(4, 26).Value = 0.8499120896328726
$ws.Cells.Item(4, 27).Value = 0.9483213126288657
$ws.Cells.Item(4, 28).Value = 0.7037095113950714
$ws.Cells.Item(4, 29).Value = 0.9060581172434008
$ws.Cells.Item(4, 30).Value = 0.7795361184671367
